# Weekly update: insert the new week's two rows (Primera / Segunda) for
# "Coliflor" at the top of the data block (row 266), pushing the existing
# rows down by two. This mirrors the source data being prepended with the
# most recent week's prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 266 and 267, shifting rows 266:359 down to 268:361
$ws.Rows("266:267").Insert()

# New row 266: Coliflor "Primera" for 2021-09-29 (serial 44468)
$ws.Cells.Item(266, 1).Value2  = 9
$ws.Cells.Item(266, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(266, 3).Value2  = "Metropolitana"
$ws.Cells.Item(266, 4).Value2  = 44468
$ws.Cells.Item(266, 5).Value2  = 13
$ws.Cells.Item(266, 6).Value2  = 100112008
$ws.Cells.Item(266, 7).Value2  = "Coliflor"
$ws.Cells.Item(266, 8).Value2  = "Sin especificar"
$ws.Cells.Item(266, 9).Value2  = "Primera"
$ws.Cells.Item(266, 10).Value2 = 2500
$ws.Cells.Item(266, 11).Value2 = 600
$ws.Cells.Item(266, 12).Value2 = 650
$ws.Cells.Item(266, 13).Value2 = 625
$ws.Cells.Item(266, 14).Value2 = "$/unidad"
$ws.Cells.Item(266, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(266, 16).Value2 = 625
$ws.Cells.Item(266, 17).Value2 = 1
$ws.Cells.Item(266, 18).Value2 = "Hortaliza"

# New row 267: Coliflor "Segunda" for 2021-09-29 (serial 44468)
$ws.Cells.Item(267, 1).Value2  = 9
$ws.Cells.Item(267, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(267, 3).Value2  = "Metropolitana"
$ws.Cells.Item(267, 4).Value2  = 44468
$ws.Cells.Item(267, 5).Value2  = 13
$ws.Cells.Item(267, 6).Value2  = 100112008
$ws.Cells.Item(267, 7).Value2  = "Coliflor"
$ws.Cells.Item(267, 8).Value2  = "Sin especificar"
$ws.Cells.Item(267, 9).Value2  = "Segunda"
$ws.Cells.Item(267, 10).Value2 = 1300
$ws.Cells.Item(267, 11).Value2 = 450
$ws.Cells.Item(267, 12).Value2 = 500
$ws.Cells.Item(267, 13).Value2 = 475
$ws.Cells.Item(267, 14).Value2 = "$/unidad"
$ws.Cells.Item(267, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(267, 16).Value2 = 475
$ws.Cells.Item(267, 17).Value2 = 1
$ws.Cells.Item(267, 18).Value2 = "Hortaliza"
